$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.100.31'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").Value = '2.338.15'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''301.27'
$ws.Range("E5").Value = '  -1.36%  '
$ws.Range("D6").Value = '''98.29'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("D7").Value = '''0.570'
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.511'
$ws.Range("E9").Value = '  -5.43%  '
$ws.Range("D10").Value = '''34.64'
$ws.Range("E10").Value = '  -3.97%  '
$ws.Range("D11").Value = '''0.0791'
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("D12").Value = '''7.10'
$ws.Range("E12").Value = '  -5.10%  '
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").Value = '2.693.04'
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = '2.331.25'
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '''13.68'
$ws.Range("E16").Value = '  -3.47%  '
$ws.Range("D17").Value = '''0.804'
$ws.Range("E17").Value = '  -4.38%  '
$ws.Range("D18").Value = '46.034.79'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = '''12.67'
$ws.Range("E19").Value = '  -6.99%  '
$ws.Range("D20").Value = '0.0₃0961'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = '''5.97'
$ws.Range("E21").Value = '  -4.08%  '
$ws.Range("D22").Value = '''66.71'
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").Value = '''245.19'
$ws.Range("E23").Value = '  -3.98%  '
$ws.Range("E24").Value = '  -5.86%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = '''1.90'
$ws.Range("E26").Value = '  -5.47%  '
$ws.Range("D27").Value = '''40.11'
$ws.Range("E27").Value = '  -4.99%  '
$ws.Range("E28").Value = '  -4.49%  '
$ws.Range("D29").Value = '''9.67'
$ws.Range("E29").Value = '  -2.99%  '
$ws.Range("D30").Value = '''20.84'
$ws.Range("E30").Value = '  +2.21%  '
$ws.Range("D31").Value = '''3.56'
$ws.Range("E31").Value = '  +13.65%  '
$ws.Range("E32").Value = '  +5.68%  '
$ws.Range("D33").Value = '''5.44'
$ws.Range("E33").Value = '  -7.43%  '
$ws.Range("D34").Value = '''144.66'
$ws.Range("E34").Value = '  -2.01%  '
$ws.Range("D35").Value = '''0.0773'
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("D36").Value = '''0.112'
$ws.Range("E36").Value = '  -2.81%  '
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("D38").Value = '''1.80'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").Value = '''15.07'
$ws.Range("E39").Value = '  +6.90%  '
$ws.Range("D40").Value = '''3.87'
$ws.Range("E40").Value = '  -3.64%  '
$ws.Range("D41").Value = '''0.0298'
$ws.Range("E41").Value = '  -4.27%  '
$ws.Range("D42").Value = '''3.20'
$ws.Range("E42").Value = '  -6.31%  '
$ws.Range("D43").Value = '''0.998'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '1.851.89'
$ws.Range("E44").Value = '  +2.83%  '
$ws.Range("D45").Value = '''90.33'
$ws.Range("E45").Value = '  -2.73%  '
$ws.Range("E46").Value = '  -8.63%  '
$ws.Range("D47").Value = '''0.185'
$ws.Range("E47").Value = '  -5.08%  '
$ws.Range("D48").Value = '''69.50'
$ws.Range("E48").Value = '  -7.13%  '
$ws.Range("D49").Value = '2.565.57'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").Value = '''96.24'
$ws.Range("E50").Value = '  -2.94%  '
$ws.Range("D51").Value = '''4.75'
$ws.Range("E51").Value = '  -2.56%  '
